$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Object Type" (column B) and "Name" (column D) columns, then
# append two new columns: Manufacturer (BattCo) and Model (Lithio 360).

for ($r = 1; $r -le 11; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 2).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $bVal
}

$ws.Range("E1").Value = "Manufacturer"
$ws.Range("F1").Value = "Model"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = "BattCo"
    $ws.Cells.Item($r, 6).Value = "Lithio 360"
}

$ws.Range("F11").Select()
